$wb = $excel.ActiveWorkbook

# Rename the "고무블럭" sheet to "고무블럭,마감재"
$ws = $wb.Worksheets.Item("고무블럭")
$ws.Name = "고무블럭,마감재"

# Make it the active sheet (moves tabSelected from 엠지스포츠 to this sheet,
# and updates workbookView's activeTab)
$ws.Activate()

# Update the selected cell on the now-active sheet to B27
[void]$ws.Range("B27").Select()
